$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Text)
    # Force the cell to store Text even when the string looks numeric
    # (plain "228.76" would otherwise be auto-coerced to a Number by Excel).
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.NumberFormat = "General"
    $Cell.Style = "Normal"
}

$ws.Range("D2").Value = '39.189.13'
$ws.Range("E2").Value = '  +1.26%  '

$ws.Range("D3").Value = '2.148.82'
$ws.Range("E3").Value = '  +2.89%  '

$ws.Range("E4").Value = '  +0.24%  '

Set-TextValue -Cell $ws.Range("D5") -Text '228.76'
$ws.Range("E5").Value = '  +0.24%  '

Set-TextValue -Cell $ws.Range("D6") -Text '0.620'
$ws.Range("E6").Value = '  +0.99%  '

Set-TextValue -Cell $ws.Range("D7") -Text '62.46'
$ws.Range("E7").Value = '  +2.94%  '

Set-TextValue -Cell $ws.Range("D8") -Text '1.00'
$ws.Range("E8").Value = '  +0.09%  '

Set-TextValue -Cell $ws.Range("D9") -Text '0.392'
$ws.Range("E9").Value = '  +1.83%  '

Set-TextValue -Cell $ws.Range("D10") -Text '0.0851'
$ws.Range("E10").Value = '  +1.80%  '

$ws.Range("E11").Value = '  -0.44%  '

Set-TextValue -Cell $ws.Range("D12") -Text '15.97'
$ws.Range("E12").Value = '  +6.89%  '

$ws.Range("D13").Value = '2.476.42'
$ws.Range("E13").Value = '  +3.15%  '

Set-TextValue -Cell $ws.Range("D14") -Text '22.23'
$ws.Range("E14").Value = '  +1.57%  '

Set-TextValue -Cell $ws.Range("D15") -Text '0.814'
$ws.Range("E15").Value = '  +2.39%  '

Set-TextValue -Cell $ws.Range("D16") -Text '5.53'
$ws.Range("E16").Value = '  +1.12%  '

$ws.Range("D17").Value = '2.154.74'
$ws.Range("E17").Value = '  +3.11%  '

$ws.Range("D18").Value = '39.215.88'
$ws.Range("E18").Value = '  +1.31%  '

$ws.Range("B19").Value = 'Litecoin'
$ws.Range("C19").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue -Cell $ws.Range("D19") -Text '71.92'
$ws.Range("E19").Value = '  +0.40%  '

$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue -Cell $ws.Range("D20") -Text '6.14'
$ws.Range("E20").Value = '  +2.01%  '

$ws.Range("D21").Value = '0.0₃0850'
$ws.Range("E21").Value = '  +1.61%  '

Set-TextValue -Cell $ws.Range("D22") -Text '227.73'
$ws.Range("E22").Value = '  +0.60%  '

$ws.Range("E23").Value = '  +0.00%  '

$ws.Range("E24").Value = '  -1.10%  '

Set-TextValue -Cell $ws.Range("D25") -Text '2.33'
$ws.Range("E25").Value = '  -0.29%  '

Set-TextValue -Cell $ws.Range("D26") -Text '9.71'
$ws.Range("E26").Value = '  +3.07%  '

Set-TextValue -Cell $ws.Range("D27") -Text '170.69'
$ws.Range("E27").Value = '  -0.03%  '

$ws.Range("E28").Value = '  -0.06%  '

Set-TextValue -Cell $ws.Range("D29") -Text '19.61'
$ws.Range("E29").Value = '  +2.39%  '

$ws.Range("E30").Value = '  -2.63%  '

Set-TextValue -Cell $ws.Range("D31") -Text '2.56'
$ws.Range("E31").Value = '  +9.70%  '

$ws.Range("E32").Value = '  +0.77%  '

Set-TextValue -Cell $ws.Range("D33") -Text '4.60'
$ws.Range("E33").Value = '  +2.23%  '

Set-TextValue -Cell $ws.Range("D34") -Text '4.82'
$ws.Range("E34").Value = '  +2.44%  '

Set-TextValue -Cell $ws.Range("D35") -Text '7.16'
$ws.Range("E35").Value = '  +12.02%  '

Set-TextValue -Cell $ws.Range("D36") -Text '0.0617'
$ws.Range("E36").Value = '  +0.62%  '

Set-TextValue -Cell $ws.Range("D37") -Text '2.41'
$ws.Range("E37").Value = '  +0.96%  '

Set-TextValue -Cell $ws.Range("D38") -Text '3.55'
$ws.Range("E38").Value = '  +0.69%  '

Set-TextValue -Cell $ws.Range("D39") -Text '0.999'
$ws.Range("E39").Value = '  -0.23%  '

Set-TextValue -Cell $ws.Range("D40") -Text '18.13'
$ws.Range("E40").Value = '  -0.65%  '

Set-TextValue -Cell $ws.Range("D41") -Text '0.0230'
$ws.Range("E41").Value = '  +3.38%  '

Set-TextValue -Cell $ws.Range("D42") -Text '102.98'
$ws.Range("E42").Value = '  +2.17%  '

$ws.Range("D43").Value = '1.535.50'
$ws.Range("E43").Value = '  -0.16%  '

Set-TextValue -Cell $ws.Range("D44") -Text '1.19'
$ws.Range("E44").Value = '  +6.44%  '

$ws.Range("B45").Value = 'HuobiToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue -Cell $ws.Range("D45") -Text '2.86'
$ws.Range("E45").Value = '  +1.59%  '

Set-TextValue -Cell $ws.Range("D46") -Text '7.84'
$ws.Range("E46").Value = '  +1.95%  '

$ws.Range("B47").Value = 'ARBITRUM'
$ws.Range("C47").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue -Cell $ws.Range("D47") -Text '1.10'
$ws.Range("E47").Value = '  +6.64%  '

$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue -Cell $ws.Range("D48") -Text '0.0917'
$ws.Range("E48").Value = '  -0.74%  '

$ws.Range("E49").Value = '  +1.08%  '

$ws.Range("D50").Value = '2.359.43'
$ws.Range("E50").Value = '  +3.07%  '

Set-TextValue -Cell $ws.Range("D51") -Text '2.96'
$ws.Range("E51").Value = '  -0.37%  '

